$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ran a 3rd data session (row 7) and entered its results ---
$ws.Range("A7").Value = 20230803
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 6
$ws.Range("G7").Value = 7

# A7 had been pre-formatted with a red "placeholder" font; now that real
# data has been entered for this session, turn the font color back to
# automatic/black so it reads like the other populated rows.
$ws.Range("A7").Font.Color = 0

# "created table" - select the new row's data as the user moves on
$ws.Range("B9").Select()

# Touch the page setup (portrait) for the sheet, as printing/layout was
# checked while wrapping up the session.
$ws.PageSetup.Orientation = 1
